# BP-359 Bank excel statemenst upload
# Reconciliation.xlsx: bank reference numbers (col B) are stored as plain
# numbers, which collide for unrelated transactions (e.g. "88888" repeated).
# They are re-entered as distinguishing text codes, and the surrounding
# data range is given explicit text / numeric formatting so Excel stops
# guessing at column types.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (REF_NO-ish bank codes): re-key as distinct text values ---
$ws.Range("B2:B6").NumberFormat = "@"
$ws.Range("B2").Value = "199999"
$ws.Range("B3").Value = "288888"
$ws.Range("B4").Value = "388888"
$ws.Range("B5").Value = "488888"
$ws.Range("B6").Value = "588888"

# --- Formatting: A:C as text, E as a fixed 2-decimal number ---
$ws.Range("A2:C6").NumberFormat = "@"
$ws.Range("E2:E6").NumberFormat = "0.00"

# --- Print orientation ---
$ws.PageSetup.Orientation = 1

# --- Active selection moves to C3 ---
$ws.Range("C3").Select() | Out-Null
